# "google - what, how, where working"
# Adds a Thursday entry to the 26/09/2022 week's existing table (row 69),
# then appends a brand-new weekly log block (header + 4 data rows) right
# below the sheet's current data, covering drive-mech / battery progress.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing "26/09/2022" week table gains a Thursday (column E) entry ---
$ws.Range("E69").Value = "working on drive mecgh"

# --- New weekly block header (row 75): date + day-of-week column titles ---
$ws.Range("A75").Value = "26/09/2022"
$ws.Range("B75").Value = "Monday"
$ws.Range("C75").Value = "Tuesday"
$ws.Range("D75").Value = "Wednesday"
$ws.Range("E75").Value = "Thursday"
$ws.Range("F75").Value = "Friday"

# --- Column B (Monday) entries, rows 76-79 ---
$ws.Range("B76").Value = "problems with cell 3"
$ws.Range("B77").Value = "tested drive mech"
$ws.Range("B78").Value = "bearing adapters are not holding"
$ws.Range("B79").Value = "try install without shaft with set screws. Also add set screws from pla to shaft to stop horizontal motion"

# --- Column C (Tuesday) entries, rows 76-77 ---
$ws.Range("C76").Value = "put drive mech together"
$ws.Range("C77").Value = "working on battery"

# --- Column D (Wednesday) entries, rows 76-79 ---
$ws.Range("D76").Value = "researching batteries"
$ws.Range("D77").Value = "going with 2 4s 3300ma and 50C"
$ws.Range("D78").Value = "designing new battery mounts "
$ws.Range("D79").Value = "created BOM"

# --- Column E (Thursday) entries, rows 76-79 ---
$ws.Range("E76").Value = "designed dome mounts"
$ws.Range("E77").Value = "working with battery"
$ws.Range("E78").Value = "found appropiate batteries"
$ws.Range("E79").Value = "worked on wed scraping for clover to answer random questions"

# --- Match the header row's look (bold/italic/shaded) to the other weekly
#     headers (row 68) by copying its formatting onto the new row 75 ---
$ws.Range("A68:F68").Copy()
$ws.Range("A75:F75").PasteSpecial(-4122)

# --- Row heights, matching the wrapped-text auto-fit heights of the new rows ---
$ws.Range("A75").RowHeight = 19.5
$ws.Range("B77").RowHeight = 30
$ws.Range("B78").RowHeight = 30
$ws.Range("B79").RowHeight = 60

# --- Leave the selection where the author left off editing ---
$ws.Range("E80").Select()
